$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.417.54"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.356.12"
$ws.Range("D3").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.629"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.62"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0920"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("E12").Value = "  -4.89%  "
$ws.Range("E13").Value = "  -4.31%  "
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.31"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.711.18"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.350.61"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.93"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +8.85%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.547.93"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.74"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +9.11%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "268.66"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.96%  "
$ws.Range("E24").Value = "  -10.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.99%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.86%  "
$ws.Range("E28").Value = "  +2.68%  "
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0902"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.37"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.09"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.59"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -8.20%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.96"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.56%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0358"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.27%  "
$ws.Range("E41").Value = "  +2.67%  "
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.22"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "117.98"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.43"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +31.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.77%  "
$ws.Range("E48").Value = "  -2.88%  "
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.568.18"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.35%  "
